$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, "A").Value = "ECs"
$ws.Cells.Item(2, "B").Value = "Il18"
$ws.Cells.Item(2, "C").Value = "Il18r1"
$ws.Cells.Item(2, "D").Value = "Inflammatory-Mac"
$ws.Cells.Item(2, "E").Value = 2
$ws.Cells.Item(2, "F").Value = 0.6666666666666666
$ws.Cells.Item(2, "G").Value = 0.420981
$ws.Cells.Item(2, "H").Value = 1.262943
$ws.Cells.Item(2, "I").Value = 0.005858863598414047
$ws.Cells.Item(2, "J").Value = 0.005858863598414048
$ws.Cells.Item(2, "K").Value = 3
$ws.Cells.Item(2, "L").Value = 1
$ws.Cells.Item(2, "M").Value = 0.499715
$ws.Cells.Item(2, "N").Value = 1.499145
$ws.Cells.Item(2, "O").Value = 0.8251882089313166
$ws.Cells.Item(2, "P").Value = 0.8251882089313167
$ws.Cells.Item(2, "Q").Value = 0.210370520415
$ws.Cells.Item(2, "R").Value = 1.893334683735
$ws.Cells.Item(2, "S").Value = 0.004834665159148177
$ws.Cells.Item(2, "T").Value = 0.004834665159148177

# Row 3
$ws.Cells.Item(3, "A").Value = "ECs"
$ws.Cells.Item(3, "B").Value = "Il18"
$ws.Cells.Item(3, "C").Value = "Il18r1"
$ws.Cells.Item(3, "D").Value = "Resolving-Mac"
$ws.Cells.Item(3, "E").Value = 2
$ws.Cells.Item(3, "F").Value = 0.6666666666666666
$ws.Cells.Item(3, "G").Value = 0.420981
$ws.Cells.Item(3, "H").Value = 1.262943
$ws.Cells.Item(3, "I").Value = 0.005858863598414047
$ws.Cells.Item(3, "J").Value = 0.005858863598414048
$ws.Cells.Item(3, "K").Value = 2
$ws.Cells.Item(3, "L").Value = 0.6666666666666666
$ws.Cells.Item(3, "M").Value = 0.105862
$ws.Cells.Item(3, "N").Value = 0.317586
$ws.Cells.Item(3, "O").Value = 0.1748117910686833
$ws.Cells.Item(3, "P").Value = 0.1748117910686833
$ws.Cells.Item(3, "Q").Value = 0.044565890622
$ws.Cells.Item(3, "R").Value = 0.401093015598
$ws.Cells.Item(3, "S").Value = 0.00102419843926587
$ws.Cells.Item(3, "T").Value = 0.00102419843926587

# Row 4
$ws.Cells.Item(4, "A").Value = "FAPs"
$ws.Cells.Item(4, "B").Value = "Il18"
$ws.Cells.Item(4, "C").Value = "Il18r1"
$ws.Cells.Item(4, "D").Value = "Inflammatory-Mac"
$ws.Cells.Item(4, "E").Value = 3
$ws.Cells.Item(4, "F").Value = 1
$ws.Cells.Item(4, "G").Value = 5.705585666666667
$ws.Cells.Item(4, "H").Value = 17.116757
$ws.Cells.Item(4, "I").Value = 0.07940559828131502
$ws.Cells.Item(4, "J").Value = 0.07940559828131502
$ws.Cells.Item(4, "K").Value = 3
$ws.Cells.Item(4, "L").Value = 1
$ws.Cells.Item(4, "M").Value = 0.499715
$ws.Cells.Item(4, "N").Value = 1.499145
$ws.Cells.Item(4, "O").Value = 0.8251882089313166
$ws.Cells.Item(4, "P").Value = 0.8251882089313167
$ws.Cells.Item(4, "Q").Value = 2.851166741418333
$ws.Cells.Item(4, "R").Value = 25.660500672765
$ws.Cells.Item(4, "S").Value = 0.06552456342487797
$ws.Cells.Item(4, "T").Value = 0.06552456342487799

# Row 5
$ws.Cells.Item(5, "A").Value = "FAPs"
$ws.Cells.Item(5, "B").Value = "Il18"
$ws.Cells.Item(5, "C").Value = "Il18r1"
$ws.Cells.Item(5, "D").Value = "Resolving-Mac"
$ws.Cells.Item(5, "E").Value = 3
$ws.Cells.Item(5, "F").Value = 1
$ws.Cells.Item(5, "G").Value = 5.705585666666667
$ws.Cells.Item(5, "H").Value = 17.116757
$ws.Cells.Item(5, "I").Value = 0.07940559828131502
$ws.Cells.Item(5, "J").Value = 0.07940559828131502
$ws.Cells.Item(5, "K").Value = 2
$ws.Cells.Item(5, "L").Value = 0.6666666666666666
$ws.Cells.Item(5, "M").Value = 0.105862
$ws.Cells.Item(5, "N").Value = 0.317586
$ws.Cells.Item(5, "O").Value = 0.1748117910686833
$ws.Cells.Item(5, "P").Value = 0.1748117910686833
$ws.Cells.Item(5, "Q").Value = 0.6040047098446668
$ws.Cells.Item(5, "R").Value = 5.436042388602001
$ws.Cells.Item(5, "S").Value = 0.01388103485643704
$ws.Cells.Item(5, "T").Value = 0.01388103485643704

# Row 6
$ws.Cells.Item(6, "A").Value = "Inflammatory-Mac"
$ws.Cells.Item(6, "B").Value = "Il18"
$ws.Cells.Item(6, "C").Value = "Il18r1"
$ws.Cells.Item(6, "D").Value = "Inflammatory-Mac"
$ws.Cells.Item(6, "E").Value = 3
$ws.Cells.Item(6, "F").Value = 1
$ws.Cells.Item(6, "G").Value = 37.57387866666667
$ws.Cells.Item(6, "H").Value = 112.721636
$ws.Cells.Item(6, "I").Value = 0.5229220082886389
$ws.Cells.Item(6, "J").Value = 0.5229220082886389
$ws.Cells.Item(6, "K").Value = 3
$ws.Cells.Item(6, "L").Value = 1
$ws.Cells.Item(6, "M").Value = 0.499715
$ws.Cells.Item(6, "N").Value = 1.499145
$ws.Cells.Item(6, "O").Value = 0.8251882089313166
$ws.Cells.Item(6, "P").Value = 0.8251882089313167
$ws.Cells.Item(6, "Q").Value = 18.77623077791333
$ws.Cells.Item(6, "R").Value = 168.98607700122
$ws.Cells.Item(6, "S").Value = 0.431509075430469
$ws.Cells.Item(6, "T").Value = 0.4315090754304691

# Row 7
$ws.Cells.Item(7, "A").Value = "Inflammatory-Mac"
$ws.Cells.Item(7, "B").Value = "Il18"
$ws.Cells.Item(7, "C").Value = "Il18r1"
$ws.Cells.Item(7, "D").Value = "Resolving-Mac"
$ws.Cells.Item(7, "E").Value = 3
$ws.Cells.Item(7, "F").Value = 1
$ws.Cells.Item(7, "G").Value = 37.57387866666667
$ws.Cells.Item(7, "H").Value = 112.721636
$ws.Cells.Item(7, "I").Value = 0.5229220082886389
$ws.Cells.Item(7, "J").Value = 0.5229220082886389
$ws.Cells.Item(7, "K").Value = 2
$ws.Cells.Item(7, "L").Value = 0.6666666666666666
$ws.Cells.Item(7, "M").Value = 0.105862
$ws.Cells.Item(7, "N").Value = 0.317586
$ws.Cells.Item(7, "O").Value = 0.1748117910686833
$ws.Cells.Item(7, "P").Value = 0.1748117910686833
$ws.Cells.Item(7, "Q").Value = 3.977645943410667
$ws.Cells.Item(7, "R").Value = 35.798813490696
$ws.Cells.Item(7, "S").Value = 0.0914129328581698
$ws.Cells.Item(7, "T").Value = 0.0914129328581698

# Row 8
$ws.Cells.Item(8, "A").Value = "MuSCs"
$ws.Cells.Item(8, "B").Value = "Il18"
$ws.Cells.Item(8, "C").Value = "Il18r1"
$ws.Cells.Item(8, "D").Value = "Inflammatory-Mac"
$ws.Cells.Item(8, "E").Value = 3
$ws.Cells.Item(8, "F").Value = 1
$ws.Cells.Item(8, "G").Value = 1.764088
$ws.Cells.Item(8, "H").Value = 5.292263999999999
$ws.Cells.Item(8, "I").Value = 0.02455111030568848
$ws.Cells.Item(8, "J").Value = 0.02455111030568848
$ws.Cells.Item(8, "K").Value = 3
$ws.Cells.Item(8, "L").Value = 1
$ws.Cells.Item(8, "M").Value = 0.499715
$ws.Cells.Item(8, "N").Value = 1.499145
$ws.Cells.Item(8, "O").Value = 0.8251882089313166
$ws.Cells.Item(8, "P").Value = 0.8251882089313167
$ws.Cells.Item(8, "Q").Value = 0.8815412349199999
$ws.Cells.Item(8, "R").Value = 7.933871114279999
$ws.Cells.Item(8, "S").Value = 0.02025928674042626
$ws.Cells.Item(8, "T").Value = 0.02025928674042626

# Row 9
$ws.Cells.Item(9, "A").Value = "MuSCs"
$ws.Cells.Item(9, "B").Value = "Il18"
$ws.Cells.Item(9, "C").Value = "Il18r1"
$ws.Cells.Item(9, "D").Value = "Resolving-Mac"
$ws.Cells.Item(9, "E").Value = 3
$ws.Cells.Item(9, "F").Value = 1
$ws.Cells.Item(9, "G").Value = 1.764088
$ws.Cells.Item(9, "H").Value = 5.292263999999999
$ws.Cells.Item(9, "I").Value = 0.02455111030568848
$ws.Cells.Item(9, "J").Value = 0.02455111030568848
$ws.Cells.Item(9, "K").Value = 2
$ws.Cells.Item(9, "L").Value = 0.6666666666666666
$ws.Cells.Item(9, "M").Value = 0.105862
$ws.Cells.Item(9, "N").Value = 0.317586
$ws.Cells.Item(9, "O").Value = 0.1748117910686833
$ws.Cells.Item(9, "P").Value = 0.1748117910686833
$ws.Cells.Item(9, "Q").Value = 0.186749883856
$ws.Cells.Item(9, "R").Value = 1.680748954704
$ws.Cells.Item(9, "S").Value = 0.00429182356526221
$ws.Cells.Item(9, "T").Value = 0.00429182356526221

# Row 10
$ws.Cells.Item(10, "A").Value = "Resolving-Mac"
$ws.Cells.Item(10, "B").Value = "Il18"
$ws.Cells.Item(10, "C").Value = "Il18r1"
$ws.Cells.Item(10, "D").Value = "Inflammatory-Mac"
$ws.Cells.Item(10, "E").Value = 3
$ws.Cells.Item(10, "F").Value = 1
$ws.Cells.Item(10, "G").Value = 26.389162
$ws.Cells.Item(10, "H").Value = 79.167486
$ws.Cells.Item(10, "I").Value = 0.3672624195259435
$ws.Cells.Item(10, "J").Value = 0.3672624195259436
$ws.Cells.Item(10, "K").Value = 3
$ws.Cells.Item(10, "L").Value = 1
$ws.Cells.Item(10, "M").Value = 0.499715
$ws.Cells.Item(10, "N").Value = 1.499145
$ws.Cells.Item(10, "O").Value = 0.8251882089313166
$ws.Cells.Item(10, "P").Value = 0.8251882089313167
$ws.Cells.Item(10, "Q").Value = 13.18706008883
$ws.Cells.Item(10, "R").Value = 118.68354079947
$ws.Cells.Item(10, "S").Value = 0.3030606181763951
$ws.Cells.Item(10, "T").Value = 0.3030606181763952

# Row 11
$ws.Cells.Item(11, "A").Value = "Resolving-Mac"
$ws.Cells.Item(11, "B").Value = "Il18"
$ws.Cells.Item(11, "C").Value = "Il18r1"
$ws.Cells.Item(11, "D").Value = "Resolving-Mac"
$ws.Cells.Item(11, "E").Value = 3
$ws.Cells.Item(11, "F").Value = 1
$ws.Cells.Item(11, "G").Value = 26.389162
$ws.Cells.Item(11, "H").Value = 79.167486
$ws.Cells.Item(11, "I").Value = 0.3672624195259435
$ws.Cells.Item(11, "J").Value = 0.3672624195259436
$ws.Cells.Item(11, "K").Value = 2
$ws.Cells.Item(11, "L").Value = 0.6666666666666666
$ws.Cells.Item(11, "M").Value = 0.105862
$ws.Cells.Item(11, "N").Value = 0.317586
$ws.Cells.Item(11, "O").Value = 0.1748117910686833
$ws.Cells.Item(11, "P").Value = 0.1748117910686833
$ws.Cells.Item(11, "Q").Value = 2.793609467644
$ws.Cells.Item(11, "R").Value = 25.142485208796
$ws.Cells.Item(11, "S").Value = 0.06420180134954834
$ws.Cells.Item(11, "T").Value = 0.06420180134954835
